$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routes detail")

$rows = @(
    @(676,40,1,614),
    @(677,40,2,615),
    @(678,40,3,616),
    @(679,40,4,617),
    @(680,40,5,618),
    @(681,40,6,619),
    @(682,40,7,620),
    @(683,40,8,621),
    @(684,40,9,622),
    @(685,40,10,623),
    @(686,108,0,614),
    @(687,109,0,623),
    @(688,41,1,624),
    @(689,41,2,625),
    @(690,41,3,626),
    @(691,41,4,627),
    @(692,41,5,628),
    @(693,42,1,629),
    @(694,42,2,630),
    @(695,42,3,631),
    @(696,42,4,632),
    @(697,42,5,633),
    @(698,42,6,634),
    @(699,42,7,635),
    @(700,42,8,636),
    @(701,42,9,637),
    @(702,42,10,638),
    @(703,42,11,639),
    @(704,42,12,640),
    @(705,42,13,641),
    @(706,42,14,642),
    @(707,42,15,643),
    @(708,42,16,644),
    @(709,110,0,624),
    @(710,110,0,625),
    @(711,110,0,626),
    @(712,110,0,627),
    @(713,110,0,645),
    @(714,110,0,646),
    @(715,110,0,647),
    @(716,110,0,648),
    @(717,110,0,649),
    @(718,110,0,650),
    @(719,110,0,651),
    @(720,110,0,652),
    @(721,43,1,652),
    @(722,43,2,653),
    @(723,43,3,654),
    @(724,43,4,655)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
}

$ws.Range("D704").Select()
